$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 96
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($null -ne $val) {
        $newVal = [int]($val / 10)
        $cell.Value2 = $newVal
    }
}

# Column A shrank (one fewer digit per value) - narrow the best-fit width
# to match the new shorter numbers.
$ws.Columns.Item(1).ColumnWidth = 4.45
